$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column D ("Tipo") to make room for "MAE"
$ws.Range("D1").EntireColumn.Insert()

# Header for new column - copy formatting from neighboring header cell (C1)
$ws.Range("D1").Value = "MAE"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# MAE values for rows 2-10
$maeValues = @(
    1.150949170179709,
    3.277320905761562,
    1.818532645477756,
    10.34970595991643,
    8.217472415074973,
    12.90282715025942,
    9.86198879267223,
    3.521111972972703,
    3.885608224300233
)

for ($i = 0; $i -lt $maeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $maeValues[$i]
}
